$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.895.96"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.304.46"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "306.53"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +2.09%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "96.65"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -1.58%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.509"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  +0.04%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.503"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -2.49%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "35.50"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  +0.19%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "18.39"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "2.662.48"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "2.312.93"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "42.858.45"
$ws.Range("E18").Value = "  -0.29%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.01"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  -1.76%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "67.38"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.84%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "236.28"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.67%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.15"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +0.06%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "4.03"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.32%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "25.34"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.29%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.30"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +12.19%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "166.70"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +1.44%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "9.08"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -0.67%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "33.16"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  +0.06%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.77"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.89%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.99"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -2.41%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "17.80"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "2.013.06"
$ws.Range("E43").Value = "  -0.44%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0280"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -2.18%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "18.15"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +3.56%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "10.07"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -3.40%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -7.04%  "
$ws.Range("E48").Value = "  -1.63%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.87"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +10.47%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "53.89"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "2.530.01"
$ws.Range("E51").Value = "  -0.01%  "
